$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'312.83"
$ws.Range('E2').Value = "'-0.52%"
$ws.Range('D3').Value = "'37.96"
$ws.Range('E3').Value = "'-3.18%"
$ws.Range('D4').Value = "'5.071"
$ws.Range('E4').Value = "'-1.46%"
$ws.Range('D5').Value = "'0.07770"
$ws.Range('E5').Value = "'-4.89%"
$ws.Range('E6').Value = "'-0.63%"
$ws.Range('D7').Value = "'1.911"
$ws.Range('E7').Value = "'-4.13%"
$ws.Range('D8').Value = "'8.197"
$ws.Range('E8').Value = "'-1.75%"
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D9').Value = "'2.980"
$ws.Range('E9').Value = "'-4.66%"
$ws.Range('B10').Value = 'MXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D10').Value = "'0.9189"
$ws.Range('E10').Value = "'-1.97%"
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = "'0.1245"
$ws.Range('E11').Value = "'-4.51%"
$ws.Range('B12').Value = 'WazirX'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D12').Value = "'0.1901"
$ws.Range('E12').Value = "'-3.57%"
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').Value = "'0.08915"
$ws.Range('E13').Value = "'-0.83%"
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').Value = "'0.03386"
$ws.Range('E14').Value = "'-4.06%"
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = "'0.09718"
$ws.Range('E15').Value = "'-0.24%"
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = "'0.001368"
$ws.Range('E16').Value = "'-2.74%"
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = "'0.006017"
$ws.Range('E17').Value = "'-8.69%"
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = "'3.539"
$ws.Range('E18').Value = "'-2.47%"
$ws.Range('D19').Value = "'0.3408"
$ws.Range('E19').Value = "'-1.76%"
$ws.Range('E20').Value = "'-0.96%"
$ws.Range('E21').Value = "'0.57%"
$ws.Range('E22').Value = "'4.03%"
$ws.Range('D23').Value = "'0.02105"
$ws.Range('E23').Value = "'5,590.24%"
$ws.Range('D24').Value = "'0.04401"
$ws.Range('E24').Value = "'0.88%"
$ws.Range('D25').Value = "'0.001214"
$ws.Range('E25').Value = "'-2.11%"
$ws.Range('D26').Value = "'0.004247"
$ws.Range('E26').Value = "'-10.75%"
$ws.Range('D27').Value = "'0.0001351"
$ws.Range('D39').Value = "'0.02139"
$ws.Range('E39').Value = "'-4.17%"
$ws.Range('D40').Value = "'0.04981"
$ws.Range('E40').Value = "'-3.87%"
$ws.Range('D41').Value = "'0.007843"
$ws.Range('E41').Value = "'1.45%"
$ws.Range('D42').Value = "'0.009893"
$ws.Range('E42').Value = "'-4.01%"
$ws.Range('D43').Value = "'0.1343"
$ws.Range('E43').Value = "'-3.98%"
$ws.Range('E44').Value = "'-1.98%"
$ws.Range('D45').Value = "'0.009673"
$ws.Range('E45').Value = "'9.02%"
$ws.Range('D46').Value = "'0.00006502"
$ws.Range('E46').Value = "'-4.71%"
$ws.Range('E47').Value = "'-0.10%"
$ws.Range('D48').Value = "'0.003073"
$ws.Range('E48').Value = "'2.11%"
$ws.Range('E49').Value = "'-0.17%"
$ws.Range('E50').Value = "'-0.10%"
$ws.Range('E51').Value = "'-0.10%"
